{"js": "// Update \"contrat pr\u00eat personne physique\" \u2014 replace ID document details\n// and the domiciliation / collection amounts (digits + French words).\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  // Replace only the first (expected unique) match, preserving its formatting.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Identity document type and number.\nawait replaceOnce(\"Passeport\", \"Carte d'identit\u00e9 nationale\");\nawait replaceOnce(\"N\u00b0PP25342A\", \"N\u00b0AA-45467776-AQ\");\n\n// Issue date and issuing authority.\nawait replaceOnce(\"14 mars 2019\", \"12 juillet 2023\");\nawait replaceOnce(\n  \"Direction g\u00e9n\u00e9rale de la documentation et l'immigation\",\n  \"Forces nationales de police\"\n);\n\n// Monthly domiciliation minimum amount (digits, then spelled-out French words).\nawait replaceOnce(\"450 000\", \"564 000\");\nawait replaceOnce(\"quatre cent cinquante mille \", \"cinq cent soixante-quatre mille \");\n\n// Collection amount for the chosen frequency (digits, then spelled-out French words).\nawait replaceOnce(\"112 500\", \"141 000\");\nawait replaceOnce(\"cent douze mille cinq cents\", \"cent quarante-et-un mille\");\n", "ps1": "# Update \"contrat pr\u00eat personne physique\" \u2014 replace ID document details\n# and the domiciliation / collection amounts (digits + French words).\n\nfunction Replace-FirstMatch {\n    param(\n        $Document,\n        [string] $SearchText,\n        [string] $NewText\n    )\n\n    $range = $Document.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $SearchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $SearchText\"\n    }\n\n    # Assign the replacement directly on the matched range (rather than via\n    # Find.Replacement + Execute(Replace:=...)) so straight apostrophes in\n    # the new text are not auto-converted into curly quotes.\n    $range.Text = $NewText\n}\n\n$d = $word.ActiveDocument\n\n# Identity document type and number.\nReplace-FirstMatch $d \"Passeport\" \"Carte d'identit\u00e9 nationale\"\nReplace-FirstMatch $d \"N\u00b0PP25342A\" \"N\u00b0AA-45467776-AQ\"\n\n# Issue date and issuing authority.\nReplace-FirstMatch $d \"14 mars 2019\" \"12 juillet 2023\"\nReplace-FirstMatch $d \"Direction g\u00e9n\u00e9rale de la documentation et l'immigation\" \"Forces nationales de police\"\n\n# Monthly domiciliation minimum amount (digits, then spelled-out French words).\nReplace-FirstMatch $d \"450 000\" \"564 000\"\nReplace-FirstMatch $d \"quatre cent cinquante mille \" \"cinq cent soixante-quatre mille \"\n\n# Collection amount for the chosen frequency (digits, then spelled-out French words).\nReplace-FirstMatch $d \"112 500\" \"141 000\"\nReplace-FirstMatch $d \"cent douze mille cinq cents\" \"cent quarante-et-un mille\"\n"}
